$wb = $excel.ActiveWorkbook

# Rename the existing sheet from "Sheet1" to "wilcox_table"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "wilcox_table"

# Add a new sheet right after it, named "wilcox_table_selected",
# containing only the rows for the selected variables.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "wilcox_table_selected"

$rows = @(
    @("Variables",    ".y.",    "group1",    "group2",    "n1", "n2", "statistic", "p"),
    @("Chla",         "Values", "OR1-1219", "OR1-1242", 5, 6, 30, 0.00433),
    @("CN",           "Values", "OR1-1219", "OR1-1242", 5, 6, 20, 0.429),
    @("D50",          "Values", "OR1-1219", "OR1-1242", 5, 6, 9,  0.329),
    @("Fluorescence", "Values", "OR1-1219", "OR1-1242", 5, 6, 18, 0.662),
    @("Porosity",     "Values", "OR1-1219", "OR1-1242", 5, 6, 6,  0.126),
    @("Temperature",  "Values", "OR1-1219", "OR1-1242", 5, 6, 10, 0.429),
    @("TOC",          "Values", "OR1-1219", "OR1-1242", 5, 6, 18, 0.662)
)

$nRows = $rows.Count
$nCols = 8
$arr = New-Object 'object[,]' $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $arr[$r,$c] = $rows[$r][$c]
    }
}

$ws2.Range("A1:H$nRows").Value = $arr

# Match the bold, centered header style used on the wilcox_table sheet.
$ws2.Range("A1:H1").Font.Bold = $true
$ws2.Range("A1:H1").HorizontalAlignment = -4108
